# Added support for Binance.
# Updates Sheet1 test rows (symbol -> BTCUSDT, date ranges, strategy rename)
# and the ListOfValues sheet (adds Binance as exchange, keeps ByBit as a
# secondary option, renames "MACD Precise" -> "Early MACD"), plus widens
# the Exchange dropdown's data validation range to include the new row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("ListOfValues")

# ---- Strategy rename first: "MACD Precise" -> "Early MACD" ----
# (applies to both test rows plus the ListOfValues lookup list)
$ws1.Range("L2").Value = "Early MACD"
$ws1.Range("L3").Value = "Early MACD"
$ws2.Range("B3").Value = "Early MACD"

# ---- New exchange: Binance (row 3 switches to it; ListOfValues row 2 becomes Binance) ----
$ws1.Range("B3").Value = "Binance"
$ws2.Range("A2").Value = "Binance"
# ByBit keeps a spot further down the ListOfValues list
$ws2.Range("A3").Value = "ByBit"

# ---- Symbol rename: BTCUSD -> BTCUSDT ----
$ws1.Range("C2").Value = "BTCUSDT"
$ws1.Range("C3").Value = "BTCUSDT"

# ---- Date range widened for both test rows ----
$ws1.Range("D2").Value = 44440
$ws1.Range("E2").Value = 44562
$ws1.Range("D3").Value = 44440
$ws1.Range("E3").Value = 44562

# ---- Widen the Exchange dropdown validation to ListOfValues!$A$2:$A$3 ----
$exchRange = $ws1.Range("B2:B1048576")
$exchRange.Validation.Delete()
$exchRange.Validation.Add(3, 1, 1, "=ListOfValues!`$A`$2:`$A`$3")
$exchRange.Validation.InputTitle = "Exchange"
$exchRange.Validation.InputMessage = "Please select exchange from dropdown list"

# ---- Selections (cosmetic, matches author's last cursor position) ----
$ws1.Range("M13").Select()
$ws2.Range("E8").Select()
